$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.7783932765807232
$ws.Range("J2").Value = 0.7783932765807231
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.572065666666667
$ws.Range("N2").Value = 4.716197
$ws.Range("O2").Value = 0.1759712293834306
$ws.Range("P2").Value = 0.1759712293834305
$ws.Range("Q2").Value = 0.3605123869426667
$ws.Range("R2").Value = 3.244611482484
$ws.Range("S2").Value = 0.1369748218237066
$ws.Range("T2").Value = 0.1369748218237065
$ws.Range("I3").Value = 0.7783932765807232
$ws.Range("J3").Value = 0.7783932765807231
$ws.Range("O3").Value = 0.4743638053196239
$ws.Range("P3").Value = 0.4743638053196239
$ws.Range("R3").Value = 8.746465288716001
$ws.Range("S3").Value = 0.3692415967140423
$ws.Range("T3").Value = 0.3692415967140423
$ws.Range("I4").Value = 0.7783932765807232
$ws.Range("J4").Value = 0.7783932765807231
$ws.Range("M4").Value = 3.123785
$ws.Range("N4").Value = 9.371355
$ws.Range("O4").Value = 0.3496649652969456
$ws.Range("P4").Value = 0.3496649652969455
$ws.Range("Q4").Value = 0.71635887134
$ws.Range("R4").Value = 6.44722984206
$ws.Range("S4").Value = 0.2721768580429744
$ws.Range("T4").Value = 0.2721768580429743
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.065288
$ws.Range("H5").Value = 0.195864
$ws.Range("I5").Value = 0.2216067234192769
$ws.Range("J5").Value = 0.2216067234192769
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.572065666666667
$ws.Range("N5").Value = 4.716197
$ws.Range("O5").Value = 0.1759712293834306
$ws.Range("P5").Value = 0.1759712293834305
$ws.Range("Q5").Value = 0.1026370232453333
$ws.Range("R5").Value = 0.9237332092080001
$ws.Range("S5").Value = 0.03899640755972403
$ws.Range("T5").Value = 0.03899640755972402
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.065288
$ws.Range("H6").Value = 0.195864
$ws.Range("I6").Value = 0.2216067234192769
$ws.Range("J6").Value = 0.2216067234192769
$ws.Range("O6").Value = 0.4743638053196239
$ws.Range("P6").Value = 0.4743638053196239
$ws.Range("Q6").Value = 0.276677551688
$ws.Range("R6").Value = 2.490097965192001
$ws.Range("S6").Value = 0.1051222086055816
$ws.Range("T6").Value = 0.1051222086055816
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.065288
$ws.Range("H7").Value = 0.195864
$ws.Range("I7").Value = 0.2216067234192769
$ws.Range("J7").Value = 0.2216067234192769
$ws.Range("M7").Value = 3.123785
$ws.Range("N7").Value = 9.371355
$ws.Range("O7").Value = 0.3496649652969456
$ws.Range("P7").Value = 0.3496649652969455
$ws.Range("Q7").Value = 0.20394567508
$ws.Range("R7").Value = 1.83551107572
$ws.Range("S7").Value = 0.07748810725397128
$ws.Range("T7").Value = 0.07748810725397125
